$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value2 = 714.7895
$ws.Range("I33").Value2 = 398.92856
$ws.Range("J33").Value2 = 1599.2
$ws.Range("K33").Value2 = 398.92856
$ws.Range("L33").Value2 = 1599.2
$ws.Range("M33").Value2 = -169.92856
$ws.Range("N33").Value2 = -2057.2
# Row 43
$ws.Range("H43").Value2 = 1036.6666
$ws.Range("I43").Value2 = 990.1667
$ws.Range("K43").Value2 = 990.1667
$ws.Range("M43").Value2 = -921.1667
# Row 70
$ws.Range("H70").Value2 = 1462882.8
$ws.Range("J70").Value2 = 6016.3335
$ws.Range("L70").Value2 = 18049.0005
$ws.Range("N70").Value2 = -18589.0005
# Row 73
$ws.Range("H73").Value2 = 1462882.8
$ws.Range("J73").Value2 = 6016.3335
$ws.Range("L73").Value2 = 18049.0005
$ws.Range("N73").Value2 = -19921.0005
# Row 92
$ws.Range("H92").Value2 = 125874.875
$ws.Range("I92").Value2 = 999
$ws.Range("K92").Value2 = 999
$ws.Range("M92").Value2 = 249

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 3929.0925
$ws.Range("I32").Value2 = 2365.4043
$ws.Range("J32").Value2 = 14428.143
$ws.Range("K32").Value2 = 2365.4043
$ws.Range("L32").Value2 = 14428.143
$ws.Range("M32").Value2 = -2078.4043
$ws.Range("N32").Value2 = -15002.143
# Row 88
$ws.Range("H88").Value2 = 27782362
$ws.Range("I88").Value2 = 41669790
$ws.Range("K88").Value2 = 41669790
$ws.Range("M88").Value2 = -41669384
# Row 91
$ws.Range("H91").Value2 = 27782362
$ws.Range("I91").Value2 = 41669790
$ws.Range("K91").Value2 = 41669790
$ws.Range("M91").Value2 = -41668386
# Row 97
$ws.Range("H97").Value2 = 512.9677
$ws.Range("I97").Value2 = 537.1539
$ws.Range("J97").Value2 = 387.2
$ws.Range("K97").Value2 = 537.1539
$ws.Range("L97").Value2 = 387.2
$ws.Range("M97").Value2 = -41.15390000000002
$ws.Range("N97").Value2 = -1379.2
# Row 122
$ws.Range("H122").Value2 = 2391.5386
$ws.Range("I122").Value2 = 1833.4706
$ws.Range("K122").Value2 = 5500.4118
$ws.Range("M122").Value2 = -3050.4118

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value2 = 7974.6553
$ws.Range("I20").Value2 = 10116.318
$ws.Range("K20").Value2 = 10116.318
$ws.Range("M20").Value2 = -9869.317999999999
# Row 40
$ws.Range("H40").Value2 = 0
$ws.Range("J40").Value2 = 0
$ws.Range("L40").Value2 = 0
$ws.Range("N40").ClearContents()
# Row 47
$ws.Range("H47").Value2 = 500000
$ws.Range("J47").Value2 = 500000
$ws.Range("L47").Value2 = 500000
$ws.Range("N47").Value2 = -501040
# Row 48
$ws.Range("H48").Value2 = 500000
$ws.Range("J48").Value2 = 500000
$ws.Range("L48").Value2 = 500000
$ws.Range("N48").Value2 = -500830
# Row 94
$ws.Range("H94").Value2 = 4947.222
$ws.Range("I94").Value2 = 7461.778
$ws.Range("J94").Value2 = 2432.6667
$ws.Range("K94").Value2 = 7461.778
$ws.Range("L94").Value2 = 2432.6667
$ws.Range("M94").Value2 = -7010.778
$ws.Range("N94").Value2 = -3334.6667
# Row 96
$ws.Range("H96").Value2 = 15346.846
$ws.Range("I96").Value2 = 15346.846
$ws.Range("K96").Value2 = 15346.846
$ws.Range("M96").Value2 = -12600.846
# Row 105
$ws.Range("H105").Value2 = 1832.16
$ws.Range("I105").Value2 = 1721.3334
$ws.Range("K105").Value2 = 1721.3334
$ws.Range("M105").Value2 = 25.66660000000002
# Row 132
$ws.Range("H132").Value2 = 96125.94
$ws.Range("J132").Value2 = 96125.94
$ws.Range("L132").Value2 = 96125.94
$ws.Range("N132").Value2 = -106245.94
# Row 134
$ws.Range("H134").Value2 = 1073.9722
$ws.Range("I134").Value2 = 1084.6571
$ws.Range("K134").Value2 = 3253.9713
$ws.Range("M134").Value2 = -718.9712999999997

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value2 = 9874.125
$ws.Range("I62").Value2 = 8427.714
$ws.Range("J62").Value2 = 19999
$ws.Range("K62").Value2 = 8427.714
$ws.Range("L62").Value2 = 19999
$ws.Range("M62").Value2 = -7803.714
$ws.Range("N62").Value2 = -21247
# Row 65
$ws.Range("H65").Value2 = 9874.125
$ws.Range("I65").Value2 = 8427.714
$ws.Range("J65").Value2 = 19999
$ws.Range("K65").Value2 = 42138.57
$ws.Range("L65").Value2 = 99995
$ws.Range("M65").Value2 = -39018.57
$ws.Range("N65").Value2 = -106235
# Row 74
$ws.Range("H74").Value2 = 57000
$ws.Range("J74").Value2 = 57000
$ws.Range("L74").Value2 = 57000
$ws.Range("N74").Value2 = -58748
# Row 77
$ws.Range("H77").Value2 = 57000
$ws.Range("J77").Value2 = 57000
$ws.Range("L77").Value2 = 171000
$ws.Range("N77").Value2 = -179736
# Row 86
$ws.Range("H86").Value2 = 44220.953
$ws.Range("I86").Value2 = 62073.734
$ws.Range("J86").Value2 = 5965
$ws.Range("K86").Value2 = 62073.734
$ws.Range("L86").Value2 = 5965
$ws.Range("M86").Value2 = -60950.734
$ws.Range("N86").Value2 = -8211
# Row 89
$ws.Range("H89").Value2 = 44220.953
$ws.Range("I89").Value2 = 62073.734
$ws.Range("J89").Value2 = 5965
$ws.Range("K89").Value2 = 310368.67
$ws.Range("L89").Value2 = 29825
$ws.Range("M89").Value2 = -304752.67
$ws.Range("N89").Value2 = -41057
# Row 99
$ws.Range("H99").Value2 = 619232.5600000001
$ws.Range("J99").Value2 = 3499
$ws.Range("L99").Value2 = 3499
$ws.Range("N99").Value2 = -6495
# Row 105
$ws.Range("H105").Value2 = 2246.2727
$ws.Range("I105").Value2 = 1663.625
$ws.Range("K105").Value2 = 1663.625
$ws.Range("M105").Value2 = 83.375
# Row 107
$ws.Range("H107").Value2 = 2215.2307
$ws.Range("I107").Value2 = 2738.6
$ws.Range("K107").Value2 = 2738.6
$ws.Range("M107").Value2 = -818.5999999999999
# Row 126
$ws.Range("H126").Value2 = 619232.5600000001
$ws.Range("J126").Value2 = 3499
$ws.Range("L126").Value2 = 10497
$ws.Range("N126").Value2 = -15437
# Row 134
$ws.Range("H134").Value2 = 720.4706
$ws.Range("I134").Value2 = 682.26666
$ws.Range("J134").Value2 = 1007
$ws.Range("K134").Value2 = 2046.79998
$ws.Range("L134").Value2 = 3021
$ws.Range("M134").Value2 = 488.20002
$ws.Range("N134").Value2 = -8091

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value2 = 838.63635
$ws.Range("I122").Value2 = 796.1667
$ws.Range("J122").Value2 = 889.6
$ws.Range("K122").Value2 = 7165.5003
$ws.Range("L122").Value2 = 8006.400000000001
$ws.Range("M122").Value2 = -4715.5003
$ws.Range("N122").Value2 = -12906.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value2 = 0
$ws.Range("J26").Value2 = 0
$ws.Range("L26").Value2 = 0
$ws.Range("N26").ClearContents()
# Row 50
$ws.Range("H50").Value2 = 0
$ws.Range("J50").Value2 = 0
$ws.Range("L50").Value2 = 0
$ws.Range("N50").ClearContents()
# Row 70
$ws.Range("H70").Value2 = 5116.4287
$ws.Range("I70").Value2 = 5192
$ws.Range("K70").Value2 = 5192
$ws.Range("M70").Value2 = -4922
# Row 73
$ws.Range("H73").Value2 = 5116.4287
$ws.Range("I73").Value2 = 5192
$ws.Range("K73").Value2 = 5192
$ws.Range("M73").Value2 = -4256
# Row 97
$ws.Range("H97").Value2 = 524.57574
$ws.Range("I97").Value2 = 517.85
$ws.Range("K97").Value2 = 517.85
$ws.Range("M97").Value2 = -21.85000000000002
# Row 113
$ws.Range("H113").Value2 = 1940
$ws.Range("I113").Value2 = 1776.7333
$ws.Range("J113").Value2 = 2348.1667
$ws.Range("K113").Value2 = 1776.7333
$ws.Range("L113").Value2 = 2348.1667
$ws.Range("M113").Value2 = 393.2666999999999
$ws.Range("N113").Value2 = -6688.1667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value2 = 29415000
$ws.Range("I7").Value2 = 50002400
$ws.Range("K7").Value2 = 50002400
$ws.Range("M7").Value2 = -50002288
# Row 16
$ws.Range("H16").Value2 = 931.3333
$ws.Range("I16").Value2 = 880.75
$ws.Range("K16").Value2 = 880.75
$ws.Range("M16").Value2 = -710.75
# Row 40
$ws.Range("H40").Value2 = 3077
$ws.Range("I40").Value2 = 2415
$ws.Range("K40").Value2 = 2415
$ws.Range("M40").Value2 = -2279
# Row 53
$ws.Range("H53").Value2 = 10994.667
$ws.Range("I53").Value2 = 7994.5
$ws.Range("K53").Value2 = 7994.5
$ws.Range("M53").Value2 = -7476.5
# Row 82
$ws.Range("H82").Value2 = 2293
$ws.Range("I82").Value2 = 2508
$ws.Range("K82").Value2 = 2508
$ws.Range("M82").Value2 = -2147
# Row 85
$ws.Range("H85").Value2 = 2293
$ws.Range("I85").Value2 = 2508
$ws.Range("K85").Value2 = 2508
$ws.Range("M85").Value2 = -1260
# Row 93
$ws.Range("H93").Value2 = 2100.3845
$ws.Range("I93").Value2 = 1550.125
$ws.Range("J93").Value2 = 2980.8
$ws.Range("K93").Value2 = 1550.125
$ws.Range("L93").Value2 = 2980.8
$ws.Range("M93").Value2 = -302.125
$ws.Range("N93").Value2 = -5476.8
# Row 122
$ws.Range("H122").Value2 = 4552.8423
$ws.Range("I122").Value2 = 3005.5
$ws.Range("J122").Value2 = 6272.1113
$ws.Range("K122").Value2 = 9016.5
$ws.Range("L122").Value2 = 18816.3339
$ws.Range("M122").Value2 = -6566.5
$ws.Range("N122").Value2 = -23716.3339
# Row 126
$ws.Range("H126").Value2 = 29415000
$ws.Range("I126").Value2 = 50002400
$ws.Range("K126").Value2 = 150007200
$ws.Range("M126").Value2 = -150004730

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value2 = 30500
$ws.Range("J41").Value2 = 30500
$ws.Range("L41").Value2 = 30500
$ws.Range("N41").Value2 = -31280
# Row 96
$ws.Range("H96").Value2 = 2700
$ws.Range("J96").Value2 = 0
$ws.Range("L96").Value2 = 0
$ws.Range("N96").ClearContents()
# Row 126
$ws.Range("H126").Value2 = 1561.4615
$ws.Range("I126").Value2 = 1474.9166
$ws.Range("K126").Value2 = 4424.7498
$ws.Range("M126").Value2 = -1954.7498
# Row 132
$ws.Range("H132").Value2 = 2388.2
$ws.Range("I132").Value2 = 2222.9395
$ws.Range("J132").Value2 = 3167.2856
$ws.Range("K132").Value2 = 6668.818499999999
$ws.Range("L132").Value2 = 9501.856800000001
$ws.Range("M132").Value2 = -4138.818499999999
$ws.Range("N132").Value2 = -14561.8568
# Row 136
$ws.Range("H136").Value2 = 1677.5714
$ws.Range("I136").Value2 = 811.17645
$ws.Range("J136").Value2 = 5359.75
$ws.Range("K136").Value2 = 2433.52935
$ws.Range("L136").Value2 = 16079.25
$ws.Range("M136").Value2 = 116.4706499999998
$ws.Range("N136").Value2 = -21179.25
